$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the three benchmarking figures in row 2 (PEPMatch row) ---
$ws.Range("B2").Value = 15.393
$ws.Range("D2").Value = 30.073
$ws.Range("E2").Value = 45.466

# --- Move the active selection to E3 (single cell), matching the saved view state ---
[void]$ws.Range("E3").Select()
